$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Product homepage paragraph: swap the two Google Code links for
#    the new GitHub links (release page + issue tracker).
# ------------------------------------------------------------------
$h1 = $d.Hyperlinks(1)
$h1.Address = "https://github.com/hardywang/batch-image-watermark-processor/releases/"
$h1.TextToDisplay = "https://github.com/hardywang/batch-image-watermark-processor/releases/"

$h2 = $d.Hyperlinks(2)
$h2.Address = "https://github.com/hardywang/batch-image-watermark-processor/issues"
$h2.TextToDisplay = "https://github.com/hardywang/batch-image-watermark-processor/issues"

$d.Content.Find.Execute("参与网上讨论。", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "与网上讨论。", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Version number paragraph: 2.12 -> 2.14 (and the File#/Build#
#    numbers that are derived from it).
# ------------------------------------------------------------------
$d.Content.Find.Execute("2. 12 (File# 2.12.0.0, Build# 2.12.4896.36977)", $true, $false, `
                         $false, $false, $false, $true, 1, $false, `
                         "2. 14 (File# 2.14.0.0, Build# 2.14.5551.18359)", 2) | Out-Null

# ------------------------------------------------------------------
# 3. System requirements paragraph: bump the required .Net Framework
#    version to 4.5, point the download link at the new Microsoft
#    page, and drop the now redundant "Service Pack 1" sentence
#    (and its hyperlink) entirely.
# ------------------------------------------------------------------
$hNet1 = $d.Hyperlinks(3)
$hNet2 = $d.Hyperlinks(4)
$spRange = $d.Range($hNet1.Range.End, $hNet2.Range.End)
$spRange.Text = ""

$d.Content.Find.Execute("请确保微软.Net Framework 3.5（下载地址", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "请确保微软.Net Framework 4.5（下载地址", 2) | Out-Null

$hNet1b = $d.Hyperlinks(3)
$hNet1b.Address = "https://www.microsoft.com/en-ca/download/details.aspx?id=30653"
$hNet1b.TextToDisplay = "https://www.microsoft.com/en-ca/download/details.aspx?id=30653"

# ------------------------------------------------------------------
# 4. Version history list: add the 2.14 release at the top.
# ------------------------------------------------------------------
$rngVH = $d.Content
$rngVH.Find.Execute("2013-05-29: 2.12") | Out-Null
$paraVH = $rngVH.Paragraphs(1)
$paraVH.Range.InsertParagraphBefore()
$rngVH2 = $d.Content
$rngVH2.Find.Execute("2013-05-29: 2.12") | Out-Null
$newVH = $rngVH2.Paragraphs(1).Previous()
$newVH.Range.Text = "2015-03-14: 2.14"

# ------------------------------------------------------------------
# 5. Change history list: document the 2.14 release (hosting move
#    from Google Code to GitHub, no functional changes).
# ------------------------------------------------------------------
$rngCH = $d.Content
$rngCH.Find.Execute("版本 2.12") | Out-Null
$paraCH = $rngCH.Paragraphs(1)
$paraCH.Range.InsertParagraphBefore()

$rngCH2 = $d.Content
$rngCH2.Find.Execute("版本 2.12") | Out-Null
$newCH1 = $rngCH2.Paragraphs(1).Previous()
$newCH1.Range.Text = "版本 2.14"

$rngCH3 = $d.Content
$rngCH3.Find.Execute("版本 2.12") | Out-Null
$paraCH3 = $rngCH3.Paragraphs(1)
$paraCH3.Range.InsertParagraphBefore()

$rngCH4 = $d.Content
$rngCH4.Find.Execute("版本 2.12") | Out-Null
$newCH2 = $rngCH4.Paragraphs(1).Previous()
$newCH2.Range.ListFormat.ListIndent()
$newCH2.Range.Text = "无功能变化，将托管服务器从Google Code移植到GitHub。"

# ------------------------------------------------------------------
# 6. Section page setup: re-assert portrait orientation (matches the
#    explicit w:orient attribute written by the editing session).
# ------------------------------------------------------------------
$d.PageSetup.Orientation = 0

